$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header E1: archivo_json -> nombre_poligono
$ws.Range("E1").Value = "nombre_poligono"

# 2) Rows 4 and 81: reclassify from No_poligono/black to Fuera/red, fill nearest polygon in I
foreach ($r in 4,81) {
    $ws.Cells.Item($r, 6).Value = "red"
    $ws.Cells.Item($r, 7).Value = "Fuera"
    $ws.Cells.Item($r, 9).Value = "29_2.json"
}

# 3) Fill column I (poligono_cercano / nombre_poligono) for remaining No_poligono rows
#    with the nearest reference polygon json file name.
$ws.Range("I7:I9").Value = "354_1.json"
$ws.Range("I10:I12").Value = "2700_1.json"
$ws.Cells.Item(14, 9).Value = "354_1.json"
$ws.Range("I16:I17").Value = "354_1.json"
$ws.Range("I18:I24").Value = "1257_3.json"
$ws.Range("I25:I28").Value = "31000_13.json"
$ws.Range("I29:I30").Value = "1257_3.json"
$ws.Range("I31:I32").Value = "2406_0.json"
$ws.Cells.Item(34, 9).Value = "2700_1.json"
$ws.Cells.Item(40, 9).Value = "354_1.json"
$ws.Range("I41:I43").Value = "2700_1.json"
$ws.Cells.Item(45, 9).Value = "2700_1.json"
$ws.Range("I67:I72").Value = "46_52.json"
$ws.Range("I77:I78").Value = "46_52.json"
$ws.Range("I84:I86").Value = "354_1.json"
$ws.Range("I87:I89").Value = "2700_1.json"
$ws.Cells.Item(92, 9).Value = "354_1.json"
$ws.Range("I94:I95").Value = "354_1.json"
$ws.Range("I96:I101").Value = "1257_3.json"
$ws.Range("I102:I106").Value = "31000_13.json"
$ws.Cells.Item(107, 9).Value = "1257_3.json"
$ws.Range("I108:I109").Value = "2406_0.json"
$ws.Cells.Item(111, 9).Value = "2700_1.json"
$ws.Cells.Item(117, 9).Value = "2406_0.json"
$ws.Cells.Item(118, 9).Value = "354_1.json"
$ws.Range("I119:I121").Value = "2700_1.json"
$ws.Range("I145:I150").Value = "46_52.json"
$ws.Range("I155:I157").Value = "46_52.json"
$ws.Cells.Item(158, 9).Value = "1361_1.json"
$ws.Cells.Item(161, 9).Value = "1361_1.json"
$ws.Range("I165:I166").Value = "847_3.json"
$ws.Cells.Item(167, 9).Value = "1361_1.json"
$ws.Cells.Item(170, 9).Value = "1361_1.json"
$ws.Range("I171:I172").Value = "847_3.json"
$ws.Range("I177:I178").Value = "847_3.json"
